# Apply "Add data for 2022-10-14" update:
# - Rename sheet title date from 2022-10-05 to 2022-10-06
# - Update header label "2022 (through 10-05)" -> "2022 (through 10-06)"
# - Update August total (I8) 164 -> 163
# - Update November total (I11) 17 -> 18

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-10-06"

$ws.Range("I1").Value = "2022 (through 10-06)"

$ws.Range("I8").Value = 163
$ws.Range("I11").Value = 18
